# Horarios actualizados Linea 141 - 734
# Applies the 08:16:28 scrape refresh to the three schedule sheets
# (LP1912, LP1912-215, 6203-6173): header timestamps/row counts, and
# the Hora_Scrap / Hora_Llegada / Linea / Minutos / Parada data rows.

$wb = $excel.ActiveWorkbook

function Set-Row {
    param($ws, $r, $a, $b, $c, $d, $e)
    $ws.Cells.Item($r, 1).Value = $a
    $ws.Cells.Item($r, 2).Value = $b
    $ws.Cells.Item($r, 3).Value = $c
    $ws.Cells.Item($r, 4).Value = $d
    $ws.Cells.Item($r, 5).Value = $e
}

# ---------------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Range("A2").Value = "Última actualización: 08:16:28"
$ws1.Range("A3").Value = "Total filas: 103"

Set-Row $ws1 74  "08:16:28" "08:21" "26_HERNANDEZ"               5   "LP1912"
Set-Row $ws1 75  "07:17:57" "08:22" "215B_EL PATO"               65  "LP1912"
Set-Row $ws1 76  "07:50:23" "08:22" "16_P MOR-SANTA ANA"         32  "LP1912"
Set-Row $ws1 77  "08:16:28" "08:23" "215B_EL PATO"               7   "LP1912"
Set-Row $ws1 79  "08:16:28" "08:27" "84_COLONIA URQUIZA-ESC 49"  11  "LP1912"
Set-Row $ws1 80  "08:16:28" "08:34" "23_HERNANDEZ"               18  "LP1912"
Set-Row $ws1 83  "08:16:28" "08:42" "81_EL PELIGRO"              26  "LP1912"
Set-Row $ws1 84  "08:16:28" "08:43" "14_ABASTO"                  27  "LP1912"
Set-Row $ws1 85  "07:50:23" "08:53" "10_OLMOS"                   63  "LP1912"
Set-Row $ws1 86  "07:17:57" "08:53" "17_ROMERO"                  96  "LP1912"
Set-Row $ws1 87  "08:16:28" "08:54" "17_ROMERO"                  38  "LP1912"
Set-Row $ws1 88  "08:16:28" "08:55" "10_OLMOS"                   39  "LP1912"
Set-Row $ws1 89  "08:16:28" "09:01" "215A_EL PATO"               45  "LP1912"
Set-Row $ws1 90  "08:16:28" "09:03" "11_ETCHEVERRY"              47  "LP1912"
Set-Row $ws1 91  "08:16:28" "09:08" "23_HERNANDEZ"               52  "LP1912"
Set-Row $ws1 92  "08:16:28" "09:10" "16_P MOR-SANTA ANA"         54  "LP1912"
Set-Row $ws1 93  "08:16:28" "09:13" "10_OLMOS"                   57  "LP1912"
Set-Row $ws1 94  "08:16:28" "09:16" "27_EL RETIRO"               60  "LP1912"
Set-Row $ws1 95  "07:50:23" "09:17" "27_EL RETIRO"               87  "LP1912"
Set-Row $ws1 96  "08:16:28" "09:21" "26_HERNANDEZ"               65  "LP1912"
Set-Row $ws1 97  "08:16:28" "09:22" "16_SANTA ANA"               66  "LP1912"
Set-Row $ws1 98  "08:16:28" "09:22" "17_ROMERO"                  66  "LP1912"
Set-Row $ws1 99  "08:16:28" "09:23" "11_ETCHEVERRY"              67  "LP1912"
Set-Row $ws1 100 "07:50:23" "09:23" "17_ROMERO"                  93  "LP1912"
Set-Row $ws1 101 "08:16:28" "09:29" "16_SANTA ANA"               73  "LP1912"
Set-Row $ws1 102 "07:50:23" "09:31" "16_SANTA ANA"               101 "LP1912"
Set-Row $ws1 103 "08:16:28" "09:32" "15_ABASTO"                  76  "LP1912"
Set-Row $ws1 104 "08:16:28" "09:33" "10_OLMOS"                   77  "LP1912"
Set-Row $ws1 105 "08:16:28" "09:42" "215C_EL PATO"               86  "LP1912"
Set-Row $ws1 106 "08:16:28" "09:43" "14_ABASTO"                  87  "LP1912"
Set-Row $ws1 107 "08:16:28" "10:10" "16_P MOR-SANTA ANA"         114 "LP1912"
Set-Row $ws1 108 "08:16:28" "10:12" "15_ABASTO"                  116 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Range("A2").Value = "Última actualización: 08:16:28"

Set-Row $ws2 23 "08:16:28" "08:23" "215B_EL PATO" 7  "LP1912"
Set-Row $ws2 24 "08:16:28" "09:01" "215A_EL PATO" 45 "LP1912"
Set-Row $ws2 25 "08:16:28" "09:42" "215C_EL PATO" 86 "LP1912"

# ---------------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Range("A2").Value = "Última actualización: 08:16:28"
$ws3.Range("A3").Value = "Total filas: 24"

Set-Row $ws3 23 "08:16:28" "08:22" "215C_LA PLATA"              6   "L6203"
Set-Row $ws3 24 "06:52:31" "08:30" "215A_LA PLATA"              98  "L6173"
Set-Row $ws3 25 "07:17:57" "08:34" "215A_LA PLATA"              77  "L6173"
Set-Row $ws3 26 "08:16:28" "08:35" "215A_LA PLATA"              19  "L6173"
Set-Row $ws3 27 "07:17:57" "09:08" "215D_LA PLATA"              111 "L6203"
Set-Row $ws3 28 "08:16:28" "09:09" "215D_LA PLATA"              53  "L6203"
Set-Row $ws3 29 "08:16:28" "10:03" "215B_LP-P MOR-40 Y 115"     107 "L6173"
